# Fix the "1=7, 2=6, 3=5, 4=4, 5=3, 6=2, 7=1" recode string: the sheet previously
# contained a duplicate shared string with a stray double space after "2=6,"
# ("1=7, 2=6,  3=5, 4=4, 5=3, 6=2, 7=1"). Retype the value in every cell that used
# it (column G "new_vals" for the umb5/umb7/umb8/umb9/umb12/umb15/umb16/umb19
# recode rows) with the corrected single-space text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedRecode = "1=7, 2=6, 3=5, 4=4, 5=3, 6=2, 7=1"

$ws.Range("G32").Value = $fixedRecode
$ws.Range("G34").Value = $fixedRecode
$ws.Range("G35").Value = $fixedRecode
$ws.Range("G36").Value = $fixedRecode
$ws.Range("G39").Value = $fixedRecode
$ws.Range("G42").Value = $fixedRecode
$ws.Range("G43").Value = $fixedRecode
$ws.Range("G46").Value = $fixedRecode

# Update the visible window state to match: scrolled down one row further and
# the active selection moved from A48 to F32.
$ws.Range("F32").Select()
